# Apply the commit: add a new "AWS" results column to the "WER por Audio"
# sheet (shifting Azure/GCP/Gemini/GPT4o/Duracao one column to the right),
# fill in previously-missing GPT4o values, convert the Duracao column from
# milliseconds to seconds, and refresh the weighted-average sheet with the
# new AWS row plus updated figures for the existing models.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet 1: "WER por Audio" ---------------------------------------------

# Insert a new column before column B; this shifts the old
# Azure/GCP/Gemini/GPT4o/Duracao columns (B:F) to (C:G).
$ws1.Range("B1").EntireColumn.Insert()

# New header for the inserted column.
$ws1.Cells.Item(1, 2).Value = "AWS"

# New AWS (%) values per row. $null means the cell stays empty, matching
# the rows that have no AWS measurement.
$awsValues = @{
    2  = "0.0%"
    3  = $null
    4  = $null
    5  = "66.67%"
    6  = "11.11%"
    7  = "42.86%"
    8  = "8.33%"
    9  = "13.59%"
    10 = "17.5%"
    11 = "28.15%"
    12 = "66.67%"
    13 = "8.98%"
    14 = "36.1%"
    15 = "2.27%"
    16 = "3.71%"
    17 = "6.29%"
    18 = "26.73%"
    19 = "12.59%"
    20 = $null
    21 = "13.33%"
    22 = "26.09%"
    23 = "65.12%"
    24 = "5.95%"
}

foreach ($row in 2..24) {
    $val = $awsValues[$row]
    if ($null -ne $val) {
        # Leading apostrophe forces the percentage-looking text to be
        # stored as literal text instead of being parsed as a number.
        $ws1.Cells.Item($row, 2).Value = "'" + $val
    }
}

# Newly-supplied GPT4o (%) values -- this column moved from E to F after
# the insert, and rows 14-24 previously had no GPT4o measurement at all.
$gpt4oValues = @{
    14 = "75.36%"
    15 = "1.42%"
    16 = "2.58%"
    17 = "5.24%"
    18 = "2.94%"
    19 = "16.08%"
    20 = "16.07%"
    21 = "1.48%"
    22 = "8.7%"
    23 = "65.12%"
    24 = "5.41%"
}

foreach ($row in 14..24) {
    $val = $gpt4oValues[$row]
    $ws1.Cells.Item($row, 6).Value = "'" + $val
}

# Duracao (s) column (now column G) -- values converted from milliseconds
# to seconds for several rows.
$durationValues = @{
    2  = 9
    3  = 5.33
    4  = 5.87
    5  = 7.42
    6  = 7.96
    7  = 5.87
    8  = 7.74
    9  = 190.2
    10 = 130.4
    11 = 93
    12 = 48.5
    13 = 174.8
    14 = 583.4
    15 = 191.8
    16 = 285.8
    17 = 253.4
    18 = 278.8
    19 = 530
    20 = 29.9
    21 = 65.59999999999999
    22 = 8
    23 = 19
    24 = 43
}

foreach ($row in 2..24) {
    $ws1.Cells.Item($row, 7).Value = $durationValues[$row]
}

# --- Sheet 2: "Media Ponderada" -------------------------------------------

# AWS now has a real weighted-average value instead of "N/A".
$ws2.Cells.Item(2, 2).Value = 18.35

# Updated weighted averages for the other models.
$ws2.Cells.Item(3, 2).Value = 15.18
$ws2.Cells.Item(4, 2).Value = 20.76
$ws2.Cells.Item(5, 2).Value = 22.47
$ws2.Cells.Item(6, 2).Value = 25.84
